$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 651.5714
$ws.Range("J17").Value = 651.5714
$ws.Range("L17").Value = 1954.7142
$ws.Range("N17").Value = -2290.7142

# Row 96: Scroll Down
$ws.Range("H96").Value = 3776.2
$ws.Range("I96").Value = 336.44446
$ws.Range("K96").Value = 1009.33338
$ws.Range("M96").Value = 363.66662

# Row 110: Make It Bigger
$ws.Range("H110").Value = 46663.332
$ws.Range("J110").Value = 59990
$ws.Range("L110").Value = 59990
$ws.Range("N110").Value = -68170

# Row 117: A Greater Grimoire
$ws.Range("H117").Value = 72499.5
$ws.Range("J117").Value = 72499.5
$ws.Range("L117").Value = 72499.5
$ws.Range("N117").Value = -81677.5

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 9717.485000000001
$ws.Range("I132").Value = 2840.8147
$ws.Range("K132").Value = 8522.444100000001
$ws.Range("M132").Value = -5992.444100000001

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 13338803
$ws.Range("I137").Value = 1499.5
$ws.Range("K137").Value = 4498.5
$ws.Range("M137").Value = -1948.5

# Row 138: All-night Crafting
$ws.Range("H138").Value = 5844.407
$ws.Range("I138").Value = 935.3333
$ws.Range("K138").Value = 2805.9999
$ws.Range("M138").Value = 2334.0001

# Row 140: Tome for Tradition
$ws.Range("H140").Value = 59748.8
$ws.Range("J140").Value = 58531
$ws.Range("L140").Value = 58531
$ws.Range("N140").Value = -68891

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 971029.1
$ws.Range("I2").Value = 1343610.1
$ws.Range("J2").Value = 2318.6
$ws.Range("K2").Value = 1343610.1
$ws.Range("L2").Value = 2318.6
$ws.Range("M2").Value = -1343497.1
$ws.Range("N2").Value = -2544.6

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7162.4873
$ws.Range("I32").Value = 4330.2964
$ws.Range("K32").Value = 4330.2964
$ws.Range("M32").Value = -4043.2964

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 11801.25
$ws.Range("I61").Value = 14235.583
$ws.Range("K61").Value = 14235.583
$ws.Range("M61").Value = -14023.583

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 8622503
$ws.Range("I74").Value = 14707164
$ws.Range("K74").Value = 14707164
$ws.Range("M74").Value = -14706290

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 8622503
$ws.Range("I77").Value = 14707164
$ws.Range("K77").Value = 73535820
$ws.Range("M77").Value = -73531452

# Row 116: No Scope
$ws.Range("H116").Value = 971029.1
$ws.Range("I116").Value = 1343610.1
$ws.Range("J116").Value = 2318.6
$ws.Range("K116").Value = 1343610.1
$ws.Range("L116").Value = 2318.6
$ws.Range("M116").Value = -1341316.1
$ws.Range("N116").Value = -6906.6

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 5304.533
$ws.Range("I122").Value = 2427.111
$ws.Range("K122").Value = 7281.333
$ws.Range("M122").Value = -4831.333

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 19192.13
$ws.Range("I132").Value = 22366.703
$ws.Range("J132").Value = 11400
$ws.Range("K132").Value = 67100.109
$ws.Range("L132").Value = 34200
$ws.Range("M132").Value = -64570.109
$ws.Range("N132").Value = -39260

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 11801.25
$ws.Range("I136").Value = 14235.583
$ws.Range("K136").Value = 42706.749
$ws.Range("M136").Value = -40156.749

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 971029.1
$ws.Range("I3").Value = 1343610.1
$ws.Range("J3").Value = 2318.6
$ws.Range("K3").Value = 1343610.1
$ws.Range("L3").Value = 2318.6
$ws.Range("M3").Value = -1343496.1
$ws.Range("N3").Value = -2546.6

# Row 5: Axe Me Anything
$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 15000
$ws.Range("K5").Value = 15000
$ws.Range("M5").Value = -14887

# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 2695.8076
$ws.Range("I20").Value = 2142.9473
$ws.Range("K20").Value = 2142.9473
$ws.Range("M20").Value = -1895.9473

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 5862.8184
$ws.Range("I134").Value = 3499
$ws.Range("K134").Value = 10497
$ws.Range("M134").Value = -7962

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2949.4
$ws.Range("I31").Value = 1186.75
$ws.Range("K31").Value = 1186.75
$ws.Range("M31").Value = -891.75

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2949.4
$ws.Range("I34").Value = 1186.75
$ws.Range("K34").Value = 1186.75
$ws.Range("M34").Value = -984.75

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 773129.5600000001
$ws.Range("I58").Value = 837140.3
$ws.Range("K58").Value = 837140.3
$ws.Range("M58").Value = -836937.3

# Row 99: O Pine
$ws.Range("H99").Value = 8583.083000000001
$ws.Range("I99").Value = 4332.6665
$ws.Range("K99").Value = 4332.6665
$ws.Range("M99").Value = -2834.6665

# Row 126: A Better Conductor
$ws.Range("H126").Value = 8583.083000000001
$ws.Range("I126").Value = 4332.6665
$ws.Range("K126").Value = 12997.9995
$ws.Range("M126").Value = -10527.9995

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 16685782
$ws.Range("I132").Value = 20848858
$ws.Range("K132").Value = 62546574
$ws.Range("M132").Value = -62544044

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3962.2
$ws.Range("I134").Value = 3976.2104
$ws.Range("K134").Value = 11928.6312
$ws.Range("M134").Value = -9393.6312

# Row 136: Turali Quality
$ws.Range("H136").Value = 773129.5600000001
$ws.Range("I136").Value = 837140.3
$ws.Range("K136").Value = 2511420.9
$ws.Range("M136").Value = -2508870.9

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 92116.62
$ws.Range("J141").Value = 96468.086
$ws.Range("L141").Value = 96468.086
$ws.Range("N141").Value = -106828.086

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 191.91667
$ws.Range("J12").Value = 185.44444
$ws.Range("L12").Value = 556.33332
$ws.Range("N12").Value = -902.33332

# Row 100: Souper
$ws.Range("H100").Value = 18666.666
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

# Row 103: West Meats East
$ws.Range("H103").Value = 1016.6667
$ws.Range("J103").Value = 2500
$ws.Range("L103").Value = 7500
$ws.Range("N103").Value = -9258

# Row 109: Cure for What Ails
$ws.Range("H109").Value = 4241.375
$ws.Range("I109").Value = 1275.8572
$ws.Range("K109").Value = 3827.5716
$ws.Range("M109").Value = -2787.5716

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 4556.4614
$ws.Range("I122").Value = 1680.5714
$ws.Range("K122").Value = 5041.7142
$ws.Range("M122").Value = -2591.7142

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 3721.3
$ws.Range("I126").Value = 3079
$ws.Range("J126").Value = 4561.231
$ws.Range("K126").Value = 9237
$ws.Range("L126").Value = 13683.693
$ws.Range("M126").Value = -6767
$ws.Range("N126").Value = -18623.693

# Row 132: On Board for Lar
$ws.Range("H132").Value = 6352
$ws.Range("I132").Value = 6568.6
$ws.Range("J132").Value = 5991
$ws.Range("K132").Value = 19705.8
$ws.Range("L132").Value = 17973
$ws.Range("M132").Value = -17175.8
$ws.Range("N132").Value = -23033

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 41669532
$ws.Range("I40").Value = 2967.6667
$ws.Range("K40").Value = 2967.6667
$ws.Range("M40").Value = -2831.6667

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 366
$ws.Range("I55").Value = 299.5
$ws.Range("K55").Value = 299.5
$ws.Range("M55").Value = -126.5

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 5709.8975
$ws.Range("I132").Value = 4284.8096
$ws.Range("K132").Value = 12854.4288
$ws.Range("M132").Value = -10324.4288

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 2852
$ws.Range("I136").Value = 2218.8572
$ws.Range("K136").Value = 6656.571599999999
$ws.Range("M136").Value = -4106.571599999999

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 4779.4243
$ws.Range("I122").Value = 4063.8518
$ws.Range("K122").Value = 12191.5554
$ws.Range("M122").Value = -9741.555399999999

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1894.25
$ws.Range("I126").Value = 1192.6666
$ws.Range("K126").Value = 3577.9998
$ws.Range("M126").Value = -1107.9998

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 39685380
$ws.Range("I132").Value = 6946819
$ws.Range("K132").Value = 20840457
$ws.Range("M132").Value = -20837927

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 9148.754000000001
$ws.Range("I136").Value = 3932.76
$ws.Range("K136").Value = 11798.28
$ws.Range("M136").Value = -9248.280000000001
